# Grade the "Driver" section (CustomerMappingDriver Class, rows 28-31):
#   row 29 -> "For successfully scanning data from input file"   -> full 16 pts
#   row 30 -> "For correct and properly aligned output"           -> full 4 pts
# The section Total (row 31) and the grand Total (row 38) are formula-driven
# and recalc automatically once the inputs below are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E29").Value = 16
$ws.Range("E30").Value = 4

# Leave the view where the grader was last working: scrolled to / selecting
# the last-graded cell in the Driver section.
$ws.Range("E30").Select()
